$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new config rows
$ws.Range("A15").Value = "pgdn_amount"
$ws.Range("B15").Value = 5

$ws.Range("A16").Value = "website_be_error"
$ws.Range("B16").Value = "We couldn't load website: {0}"

# Update selection to match the post-edit state
$ws.Range("L25").Select()
